# Update "想去人数" (F column, want-to-go counts) across all four sheets and
# refresh the newly-scraped 萤火虫动漫游戏嘉年华 x KKWORLD2024 event row in the
# "全部类型" sheet: it now sorts in before 冰兔2024/跨越二次元ACG (by start date),
# which pushes those two rows down one slot and drops the old row 15
# (音波狂潮II 萤光宇宙 音游嘉年华) — matches "Update gh-pages to output generated
# at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1489
$ws.Range("F3").Value = 1460
$ws.Range("F5").Value = 230
$ws.Range("F6").Value = 726
$ws.Range("F7").Value = 41
$ws.Range("F8").Value = 653
$ws.Range("F9").Value = 472
$ws.Range("F11").Value = 1389
$ws.Range("F12").Value = 36242
$ws.Range("F13").Value = 7202
$ws.Range("F14").Value = 117
$ws.Range("F15").Value = 378
$ws.Range("F16").Value = 587
$ws.Range("F17").Value = 451
$ws.Range("F19").Value = 116
$ws.Range("F20").Value = 256
$ws.Range("F21").Value = 52
$ws.Range("F22").Value = 457
$ws.Range("F23").Value = 117
$ws.Range("F24").Value = 820
$ws.Range("F25").Value = 20
$ws.Range("F26").Value = 324
$ws.Range("F27").Value = 396
$ws.Range("F28").Value = 450
$ws.Range("F29").Value = 25
$ws.Range("F30").Value = 222
$ws.Range("F31").Value = 57
$ws.Range("F32").Value = 748
$ws.Range("F33").Value = 294
$ws.Range("F34").Value = 136
$ws.Range("F35").Value = 762
$ws.Range("F36").Value = 116
$ws.Range("F38").Value = 808
$ws.Range("F39").Value = 295
$ws.Range("F40").Value = 54
$ws.Range("F41").Value = 27

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1215
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 169
$ws.Range("F7").Value = 4330
$ws.Range("F9").Value = 243
$ws.Range("F13").Value = 41
$ws.Range("F19").Value = 4304

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1478
$ws.Range("F3").Value = 364

$ws = $wb.Worksheets.Item("全部类型")

# New event inserted at row 13 (萤火虫动漫游戏嘉年华 x KKWORLD2024); the two rows
# that used to occupy 13-14 (冰兔2024 / 跨越二次元ACG) shift down to 14-15, and the
# former row 15 (音波狂潮II) drops off the sheet. Column A (row index) is untouched.
$ws.Range("B13").Value = "'2024-07-19"
$ws.Range("C13").Value = "广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园"
$ws.Range("D13").Value = "新港东路1000号 保利世贸博览馆"
$ws.Range("E13").Value = "2024.07.19 09:00-07.22 17:00"
$ws.Range("F13").Value = 36244
$ws.Range("G13").Value = "暂时售罄"
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=87210"
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg"

$ws.Range("B14").Value = "'2024-07-20"
$ws.Range("C14").Value = "广州·冰兔2024线下live「过去和未来」"
$ws.Range("D14").Value = "恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）"
$ws.Range("E14").Value = "2024.07.20 20:00-07.20 22:00"
$ws.Range("F14").Value = 169
$ws.Range("G14").Value = 198
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=87546"
$ws.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg"

$ws.Range("B15").Value = "'2024-07-20"
$ws.Range("C15").Value = "广州·跨越二次元ACG神级动漫世界巡回演唱会"
$ws.Range("D15").Value = "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
$ws.Range("E15").Value = "2024.07.20 19:30-07.20 21:10"
$ws.Range("F15").Value = 292
$ws.Range("G15").Value = 280
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=85353"
$ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"

# Remaining want-to-go count refreshes for rows untouched by the shift above.
$ws.Range("F2").Value = 1478
$ws.Range("F3").Value = 364
$ws.Range("F4").Value = 1215
$ws.Range("F5").Value = 1489
$ws.Range("F7").Value = 1460
$ws.Range("F8").Value = 230
$ws.Range("F9").Value = 726
$ws.Range("F10").Value = 41
$ws.Range("F11").Value = 653
$ws.Range("F12").Value = 2
$ws.Range("F16").Value = 243
$ws.Range("F17").Value = 243
$ws.Range("F20").Value = 7202
$ws.Range("F21").Value = 378
$ws.Range("F23").Value = 587
$ws.Range("F24").Value = 451
$ws.Range("F25").Value = 41
$ws.Range("F26").Value = 116
$ws.Range("F27").Value = 256
$ws.Range("F29").Value = 52
$ws.Range("F31").Value = 457
$ws.Range("F32").Value = 117
$ws.Range("F33").Value = 820
$ws.Range("F34").Value = 20
$ws.Range("F35").Value = 324
$ws.Range("F36").Value = 396
$ws.Range("F37").Value = 450
$ws.Range("F38").Value = 25
$ws.Range("F39").Value = 222
$ws.Range("F40").Value = 57
$ws.Range("F41").Value = 748
$ws.Range("F43").Value = 294
$ws.Range("F44").Value = 136
$ws.Range("F45").Value = 808
$ws.Range("F46").Value = 295
$ws.Range("F47").Value = 54
$ws.Range("F49").Value = 27
